# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.505614041169197
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 16.98373111632243
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 267.1278323104921
